# person.xlsx -- "Prepared to commit as 'new' maestro repo in GitHub
# maestrobpm organization"
#
# Semantic changes captured from the OOXML diff:
#   1. The third worksheet is renamed from "res.usrs" -> "res.users"
#      (a typo fix picked up while re-publishing the workbook).
#   2. The workbook's active/selected tab moves from "res.users" (the
#      sheet that was selected when the file was last saved) back to
#      "Master", the first sheet -- i.e. <sheetView tabSelected="1"/>
#      moves from sheet3 to sheet1 and workbookView's activeTab reverts
#      to the first sheet.
#
# (The remaining hunks in the diff -- windowHeight, per-column width
# jitter of a few hundredths of a character, and the x14ac:dyDescent /
# default-row-height bump on every sheet -- are Excel's own pixel-metric
# re-layout of the grid, not a user edit; they are recomputed by the
# Excel engine from live font metrics and aren't exposed as settable
# properties on the object model, so there is nothing to replay here.)

$wb = $excel.ActiveWorkbook

# 1. Fix the worksheet name.
$wsUsers = $wb.Worksheets.Item("res.usrs")
$wsUsers.Name = "res.users"

# 2. Re-select "Master" as the active sheet/tab.
$wsMaster = $wb.Worksheets.Item("Master")
$wsMaster.Activate()
